$p = $ppt.ActivePresentation

# --- Slide 11: rewrite the task body to the merged/expanded wording ---
$s11 = $p.Slides.Item(11)
$content = $s11.Shapes.Item(2).TextFrame.TextRange

$nl = [char]13
$text  = "Parašykite procedūras žemiau esančioms užklausoms. Remkites matytas pavyzdžiais." + $nl
$text += "1. Kiek yra pardavinėjama skirtingų produktų? " + $nl
$text += "2. Atraskite klientus, kurie neturi pardavimų atstovo. Grąžinkite klientų vardus ir miestus iš kur jie yra. " + $nl
$text += "3. Kokie yra VP ir Managers vardai? Vardus ir pavardes išveskite viename stulpelyje." + $nl
$text += "4. Raskite klientus, kurie yra iš nurodyto miesto. " + $nl
$text += "Pvz. call ClientsByCity('Frankfurt');" + $nl

$content.Text = $text

# --- Slide 13 ("Užduotis" duplicate with the old wording) is folded into
#     slide 11 above, so drop the now-redundant slide. ---
$s13 = $p.Slides.Item(13)
$s13.Delete()
